$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 14398.54268594074
$ws.Range("C2").Value = 19229.14446931734
$ws.Range("D2").Value = 2206.837529930249
